# Trade #11 closed at 2026-02-17 15:14:36 - unknown UNKNOWN +0.000%
#
# Updates Summary + Strategy Status roll-up figures and appends the new
# closed trade (Trade #11 / row 12, 0-based trade id 11) to both the
# "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# --- Summary sheet --------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.24    # Total P&L %
$summary.Range("B6").Value = 11       # Total Trades
$summary.Range("B9").Value = 27.27    # Win Rate %

# --- Strategy Status sheet ------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 11        # MarketMaking Trades
$status.Range("G4").Value = 27.27     # MarketMaking Win Rate %

# --- New trade row appended to "All Trades" and "MarketMaking" ------
$newRow = @(11, "2026-02-17", "15:14:30", "MarketMaking", "DOWN", 0.02, 0.02, "CLOSED", 0, 0, 99.87, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    # The date column (B) looks like a date string ("2026-02-17") and Excel
    # would otherwise auto-convert it to a date serial number on assignment;
    # force the cell to Text format first so it is stored as a literal string,
    # matching the rest of the column.
    $ws.Cells.Item(12, 2).NumberFormat = "@"
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $ws.Cells.Item(12, $i + 1).Value = $newRow[$i]
    }
}
